# "all pricings and features are now read from csv"
#
# This script:
#  1. On the "ondemand" sheet, inserts a new "currency" column (with value
#     "USD" for every data row) right after the "Date" column.
#  2. On the "features" and "features tocsv" sheets, renames the first
#     header column from "Keyword" to "apiname" (values already held the
#     api name, only the header text was stale).
#  3. Restores the view state (active sheet / selected cells) to match the
#     edited workbook.

$wb = $excel.ActiveWorkbook

$ondemand = $wb.Worksheets.Item(1)   # "ondemand"
$features = $wb.Worksheets.Item(2)   # "features"
$featuresToCsv = $wb.Worksheets.Item(3)   # "features tocsv"

# --- 1. "ondemand": insert a "currency" column (USD) after column A (Date) ---
$ondemand.Columns("B:B").Insert()
$ondemand.Columns("B:B").ColumnWidth = $ondemand.Columns("A:A").ColumnWidth

$ondemand.Range("B1").Value2 = "currency"
$ondemand.Range("B2:B5").Value2 = "USD"

# --- 2. "features" / "features tocsv": header rename Keyword -> apiname ---
$features.Range("A1").Value2 = "apiname"
$featuresToCsv.Range("A1").Value2 = "apiname"

# --- 3. view / selection state ---
$featuresToCsv.Range("A1").Select()

$features.Activate()
$features.Range("A2").Select()

$ondemand.Activate()
$ondemand.Range("G21").Select()
